$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised cumulative AgTests (F) / AgPosit (G) figures for existing rows ---
$ws.Range("F630").Value = 47067
$ws.Range("G630").Value = 2992
$ws.Range("F631").Value = 42267
$ws.Range("G631").Value = 2781
$ws.Range("F656").Value = 52690
$ws.Range("G656").Value = 1244
$ws.Range("F684").Value = 57419
$ws.Range("G684").Value = 1218
$ws.Range("F691").Value = 62647
$ws.Range("F711").Value = 22596
$ws.Range("G711").Value = 3821
$ws.Range("F722").Value = 27946
$ws.Range("F723").Value = 22548
$ws.Range("G723").Value = 2750
$ws.Range("F725").Value = 12753
$ws.Range("G725").Value = 2087
$ws.Range("F726").Value = 35804
$ws.Range("G726").Value = 4104
$ws.Range("F727").Value = 25087
$ws.Range("G727").Value = 2803
$ws.Range("F728").Value = 24631
$ws.Range("G728").Value = 2608
$ws.Range("F729").Value = 23223
$ws.Range("G729").Value = 2516
$ws.Range("F730").Value = 19460
$ws.Range("G730").Value = 2324
$ws.Range("F731").Value = 8609
$ws.Range("G731").Value = 1322
$ws.Range("F732").Value = 11827
$ws.Range("G732").Value = 1896
$ws.Range("F733").Value = 31397
$ws.Range("G733").Value = 3690
$ws.Range("F734").Value = 23039
$ws.Range("G734").Value = 2527
$ws.Range("F735").Value = 18618
$ws.Range("G735").Value = 2182
$ws.Range("F736").Value = 18359
$ws.Range("G736").Value = 2068

# --- New rows 737-739: daily stats through 2022-03-13 (po 14. 03. 2022 update) ---
$ws.Range("A737").Value = 44631
$ws.Range("B737").Value = 1565484
$ws.Range("C737").Value = 16319
$ws.Range("D737").Value = 9051
$ws.Range("E737").Value = 18881
$ws.Range("F737").Value = 14158
$ws.Range("G737").Value = 1960

$ws.Range("A738").Value = 44632
$ws.Range("B738").Value = 1573111
$ws.Range("C738").Value = 13475
$ws.Range("D738").Value = 7627
$ws.Range("E738").Value = 18918
$ws.Range("F738").Value = 4157
$ws.Range("G738").Value = 716

$ws.Range("A739").Value = 44633
$ws.Range("B739").Value = 1576486
$ws.Range("C739").Value = 6307
$ws.Range("D739").Value = 3375
$ws.Range("E739").Value = 18950
$ws.Range("F739").Value = 4311
$ws.Range("G739").Value = 880
